{"js": "// Replace the \"observation campaign dates\" paragraph (and its trailing\n// constellation sentence) with the new Gemini campaign text. This occurs\n// 4 times in the document; each time the paragraph's multiple runs are\n// collapsed into a single, formatting-less run.\n\nconst searchText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od\";\nconst newText = \"Informace v t\u00e9to p\u0159\u00edru\u010dce jsou ur\u010deny pro pozorovac\u00ed kampa\u0148 prob\u00edhaj\u00edc\u00ed od Gemini: 14.-23. \u00fanora, 14.-24. b\u0159ezna\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\n// Resolve the full paragraph for every match first (search ranges become\n// unusable after the underlying content is mutated).\nconst paragraphs = [];\nfor (let i = 0; i < results.items.length; i++) {\n  const para = results.items[i].paragraphs.getFirst();\n  para.load(\"text\");\n  paragraphs.push(para);\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.length; i++) {\n  const para = paragraphs[i];\n  // Only touch paragraphs that actually start with the target sentence,\n  // to avoid accidentally affecting unrelated matches.\n  if (para.text.indexOf(searchText) === 0) {\n    para.clear();\n    para.insertText(newText, \"Start\");\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the \"observation campaign dates\" paragraph (and its trailing\n# constellation sentence) with the new Gemini campaign text. This occurs\n# 4 times in the document; each time the paragraph's multiple runs are\n# collapsed into a single, formatting-less run.\n\n$d = $word.ActiveDocument\n\n$searchText = \"Informace v t\"+[char]0xE9+\"to p\"+[char]0x159+[char]0xED+\"ru\"+[char]0x10D+\"ce jsou ur\"+[char]0x10D+\"eny pro pozorovac\"+[char]0xED+\" kampa\"+[char]0x148+\" prob\"+[char]0xED+\"haj\"+[char]0xED+\"c\"+[char]0xED+\" od\"\n$newText = \"Informace v t\"+[char]0xE9+\"to p\"+[char]0x159+[char]0xED+\"ru\"+[char]0x10D+\"ce jsou ur\"+[char]0x10D+\"eny pro pozorovac\"+[char]0xED+\" kampa\"+[char]0x148+\" prob\"+[char]0xED+\"haj\"+[char]0xED+\"c\"+[char]0xED+\" od Gemini: 14.-23. \"+[char]0xFA+\"nora, 14.-24. b\"+[char]0x159+\"ezna\"\n\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith($searchText)) {\n        $targets += $p\n    }\n}\n\nforeach ($p in $targets) {\n    $r = $p.Range\n    # Exclude the trailing paragraph mark from the range.\n    $r.End = $r.End - 1\n    $r.Delete()\n\n    $r2 = $p.Range\n    $r2.End = $r2.End - 1\n    $r2.InsertAfter($newText)\n}\n"}
